$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("Q2").Value = 1.55
$ws.Range("BH2").Value = "2026-02-24 14:29:51"

# Row 3
$ws.Range("G3").Value = 2.1
$ws.Range("BH3").Value = "2026-02-24 14:29:51"

# Row 4
$ws.Range("BH4").Value = "2026-02-24 14:29:51"

# Row 5
$ws.Range("Q5").Value = 2.16
$ws.Range("BH5").Value = "2026-02-24 14:29:51"

# Row 6
$ws.Range("BH6").Value = "2026-02-24 14:29:51"

# Row 7
$ws.Range("I7").Value = 4.1
$ws.Range("J7").Value = 3.4
$ws.Range("K7").Value = 3.65
$ws.Range("BH7").Value = "2026-02-24 14:29:51"

# Row 8
$ws.Range("Q8").Value = 3.5
$ws.Range("BH8").Value = "2026-02-24 14:29:51"

# Row 9
$ws.Range("F9").Value = 1.5
$ws.Range("BH9").Value = "2026-02-24 14:29:51"

# Row 10
$ws.Range("BH10").Value = "2026-02-24 14:29:51"
